$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting existing D:K data to E:L
$ws.Range("D:D").Insert()

# Copy number formats/styles from column E into the newly inserted column D
$ws.Range("E:E").Copy()
$ws.Range("D:D").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D with the latest fiscal-year figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 4043700
$ws.Range("D9").Value = 491000
$ws.Range("D10").Value = 3552700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 230400
$ws.Range("D17").Value = 3226600
$ws.Range("D18").Value = 817100
$ws.Range("D20").Value = -1300
$ws.Range("D21").Value = 1046100
$ws.Range("D22").Value = 200
$ws.Range("D23").Value = 815500
$ws.Range("D24").Value = 209800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 605700
$ws.Range("D27").Value = 605700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 1300
$ws.Range("D33").Value = 605700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 605700
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 190300
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 468300
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 47700
$ws.Range("D46").Value = 706200
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 2754900
$ws.Range("D49").Value = 19500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 64600
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 3545300
$ws.Range("D57").Value = 78500
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 278200
$ws.Range("D60").Value = 356700
$ws.Range("D61").Value = 45000
$ws.Range("D62").Value = 463100
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 864800
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 2530200
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2680500
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 605700
$ws.Range("D83").Value = 230400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 900100
$ws.Range("D91").Value = -588300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -580400
$ws.Range("D96").Value = -42600
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -256900
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 62800
